# Update the home advantage attribute
# Set the Date column (A2:A19) to the same date (2025-08-22) for all matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = Get-Date -Year 2025 -Month 8 -Day 22 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 1).Value = $newDate
}
